$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (Volume/Number + Report week dates) ----
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# ---- Simple numeric value updates (same type before/after) ----
$ws.Range("D15").Value = 2
$ws.Range("G15").Value = 7
$ws.Range("J15").Value = 5
$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 140
$ws.Range("F16").Value = 48
$ws.Range("G16").Value = 37
$ws.Range("H16").Value = 29.729729729729
$ws.Range("I16").Value = 23
$ws.Range("J16").Value = 16
$ws.Range("K16").Value = 43.75
$ws.Range("L16").Value = 91.666666666666
$ws.Range("M16").Value = 35.294117647058
$ws.Range("N16").Value = -72.619047619047
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -35.294117647058
$ws.Range("F17").Value = 58
$ws.Range("H17").Value = -13.432835820895
$ws.Range("I17").Value = 24
$ws.Range("J17").Value = 36
$ws.Range("K17").Value = -33.333333333333
$ws.Range("L17").Value = 20
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = -35.135135135135
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 31
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = 29.166666666666
$ws.Range("I18").Value = 13
$ws.Range("J18").Value = 11
$ws.Range("K18").Value = 18.181818181818
$ws.Range("L18").Value = 18.181818181818
$ws.Range("M18").Value = 44.444444444444
$ws.Range("N18").Value = -77.966101694915
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 71.428571428571
$ws.Range("F19").Value = 75
$ws.Range("H19").Value = 17.1875
$ws.Range("I19").Value = 43
$ws.Range("J19").Value = 30
$ws.Range("K19").Value = 43.333333333333
$ws.Range("L19").Value = 48.275862068965
$ws.Range("M19").Value = 186.666666666667
$ws.Range("N19").Value = 65.384615384615
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 20
$ws.Range("G20").Value = 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 9
$ws.Range("J20").Value = 8
$ws.Range("K20").Value = 12.5
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 80
$ws.Range("N20").Value = -66.666666666666
$ws.Range("C21").Value = 58
$ws.Range("D21").Value = 51
$ws.Range("E21").Value = 13.725490196078
$ws.Range("F21").Value = 232
$ws.Range("G21").Value = 219
$ws.Range("H21").Value = 5.936073059360
$ws.Range("I21").Value = 112
$ws.Range("J21").Value = 106
$ws.Range("K21").Value = 5.660377358490
$ws.Range("L21").Value = 36.585365853658
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = -53.138075313807
$ws.Range("C22").Value = 3
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 4
$ws.Range("L22").Value = -33.333333333333
$ws.Range("M22").Value = 100
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 30
$ws.Range("G23").Value = 36
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 13
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = -23.529411764705
$ws.Range("L23").Value = 18.181818181818
$ws.Range("M23").Value = 85.714285714285
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = -9.756097560975
$ws.Range("F24").Value = 149
$ws.Range("G24").Value = 142
$ws.Range("H24").Value = 4.929577464788
$ws.Range("I24").Value = 75
$ws.Range("J24").Value = 65
$ws.Range("K24").Value = 15.384615384615
$ws.Range("L24").Value = 82.926829268292
$ws.Range("M24").Value = 87.5
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 27.272727272727
$ws.Range("F25").Value = 91
$ws.Range("G25").Value = 72
$ws.Range("H25").Value = 26.388888888888
$ws.Range("I25").Value = 40
$ws.Range("J25").Value = 26
$ws.Range("K25").Value = 53.846153846153
$ws.Range("L25").Value = 60
$ws.Range("M25").Value = 17.647058823529
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 9
$ws.Range("H26").Value = -88.888888888888
$ws.Range("J26").Value = 5
$ws.Range("K26").Value = -80
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -20
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 3
$ws.Range("L27").Value = 50
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -66.666666666666
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J38").Value = 668
$ws.Range("K38").Value = 15.771230502599
$ws.Range("L38").Value = -15.762925598991
$ws.Range("M38").Value = -64.989517819706
$ws.Range("N38").Value = -69.102682701202
$ws.Range("J40").Value = 326
$ws.Range("K40").Value = 2.839116719242
$ws.Range("L40").Value = -38.490566037735
$ws.Range("M40").Value = -76.099706744868
$ws.Range("N40").Value = -76.27365356623

# ---- Cells changing from text placeholder ("***.*") to numeric percent value ----
# Reuses the workbook's existing percent number format so the style index (15) is preserved.
$ws.Range("L15").Value = -100
$ws.Range("L15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L26").Value = 0
$ws.Range("L26").NumberFormat = '#,##0.0;"-"#,##0.0'

# ---- Cells changing from numeric value to text placeholder ("0" or "***.*") ----
# Force a text-typed cell via the text number format, then assign the literal string.
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
